# Swap the contents of columns C (codeforiati:group-code) and D
# (codeforiati:group-name) for every row in the sheet, including the header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value()
    $dVal = $dCell.Value()

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
